# Swap the Id (A), Ost (Q) and Nord (R) values between row 2 and row 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$cols = @("A", "Q", "R")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $val2 = $cell2.Value2
    $val3 = $cell3.Value2

    $cell2.Value2 = $val3
    $cell3.Value2 = $val2
}
